$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1697.8889
$ws.Range("I107").Value = 276.5
$ws.Range("K107").Value = 276.5
$ws.Range("M107").Value = 1643.5

$ws.Range("H113").Value = 55566890
$ws.Range("I113").Value = 16672165
$ws.Range("J113").Value = 133356340
$ws.Range("K113").Value = 16672165
$ws.Range("L113").Value = 133356340
$ws.Range("M113").Value = -16668911
$ws.Range("N113").Value = -133362848

$ws.Range("H132").Value = 2498.2222
$ws.Range("I132").Value = 2216.7812
$ws.Range("K132").Value = 6650.3436
$ws.Range("M132").Value = -4120.3436

$ws.Range("H137").Value = 3085.8484
$ws.Range("I137").Value = 1282.5555
$ws.Range("K137").Value = 3847.6665
$ws.Range("M137").Value = -1297.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14292511
$ws.Range("I32").Value = 14292511
$ws.Range("K32").Value = 14292511
$ws.Range("M32").Value = -14292224

$ws.Range("H45").Value = 2203
$ws.Range("I45").Value = 2253.5
$ws.Range("K45").Value = 2253.5
$ws.Range("M45").Value = -1876.5

$ws.Range("H61").Value = 21790330
$ws.Range("I61").Value = 62503388
$ws.Range("K61").Value = 62503388
$ws.Range("M61").Value = -62503176

$ws.Range("H74").Value = 13900107
$ws.Range("I74").Value = 35715284
$ws.Range("K74").Value = 35715284
$ws.Range("M74").Value = -35714410

$ws.Range("H77").Value = 13900107
$ws.Range("I77").Value = 35715284
$ws.Range("K77").Value = 178576420
$ws.Range("M77").Value = -178572052

$ws.Range("H136").Value = 21790330
$ws.Range("I136").Value = 62503388
$ws.Range("K136").Value = 187510164
$ws.Range("M136").Value = -187507614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 896.0323
$ws.Range("I94").Value = 896.0323
$ws.Range("K94").Value = 896.0323
$ws.Range("M94").Value = -445.0323

$ws.Range("H105").Value = 2263.2222
$ws.Range("I105").Value = 1714.25
$ws.Range("J105").Value = 2702.4
$ws.Range("K105").Value = 1714.25
$ws.Range("L105").Value = 2702.4
$ws.Range("M105").Value = 32.75
$ws.Range("N105").Value = -6196.4

$ws.Range("H134").Value = 33832
$ws.Range("I134").Value = 1507.52
$ws.Range("K134").Value = 4522.559999999999
$ws.Range("M134").Value = -1987.559999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 632563.8
$ws.Range("I31").Value = 1843.3334
$ws.Range("J31").Value = 1460384.5
$ws.Range("K31").Value = 1843.3334
$ws.Range("L31").Value = 1460384.5
$ws.Range("M31").Value = -1548.3334
$ws.Range("N31").Value = -1460974.5

$ws.Range("H34").Value = 632563.8
$ws.Range("I34").Value = 1843.3334
$ws.Range("J34").Value = 1460384.5
$ws.Range("K34").Value = 1843.3334
$ws.Range("L34").Value = 1460384.5
$ws.Range("M34").Value = -1641.3334
$ws.Range("N34").Value = -1460788.5

$ws.Range("H58").Value = 2422.6667
$ws.Range("I58").Value = 2291.4583
$ws.Range("K58").Value = 2291.4583
$ws.Range("M58").Value = -2088.4583

$ws.Range("H105").Value = 1997.1111
$ws.Range("I105").Value = 1727.25
$ws.Range("K105").Value = 1727.25
$ws.Range("M105").Value = 19.75

$ws.Range("H106").Value = 28280
$ws.Range("J106").Value = 28280
$ws.Range("L106").Value = 28280
$ws.Range("N106").Value = -30804

$ws.Range("H107").Value = 3002.4
$ws.Range("I107").Value = 2999.75
$ws.Range("K107").Value = 2999.75
$ws.Range("M107").Value = -1079.75

$ws.Range("H122").Value = 3412.1667
$ws.Range("J122").Value = 3057
$ws.Range("L122").Value = 9171
$ws.Range("N122").Value = -14071

$ws.Range("H132").Value = 3113.6667
$ws.Range("I132").Value = 2752.875
$ws.Range("K132").Value = 8258.625
$ws.Range("M132").Value = -5728.625

$ws.Range("H136").Value = 2422.6667
$ws.Range("I136").Value = 2291.4583
$ws.Range("K136").Value = 6874.374899999999
$ws.Range("M136").Value = -4324.374899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 528907.25
$ws.Range("J12").Value = 792350.4399999999
$ws.Range("L12").Value = 2377051.32
$ws.Range("N12").Value = -2377397.32

$ws.Range("H39").Value = 642104.7
$ws.Range("I39").Value = 426315
$ws.Range("J39").Value = 749999.5
$ws.Range("K39").Value = 1278945
$ws.Range("L39").Value = 2249998.5
$ws.Range("M39").Value = -1278651
$ws.Range("N39").Value = -2250586.5

$ws.Range("H80").Value = 3960.25
$ws.Range("I80").Value = 2993.3333
$ws.Range("K80").Value = 8979.999899999999
$ws.Range("M80").Value = -8043.999899999999

$ws.Range("H82").Value = 5821.25
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H83").Value = 3960.25
$ws.Range("I83").Value = 2993.3333
$ws.Range("K83").Value = 26939.9997
$ws.Range("M83").Value = -22259.9997

$ws.Range("H85").Value = 5821.25
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1528.5
$ws.Range("I107").Value = 1473.875
$ws.Range("J107").Value = 1747
$ws.Range("K107").Value = 1473.875
$ws.Range("L107").Value = 1747
$ws.Range("M107").Value = 446.125
$ws.Range("N107").Value = -5587

$ws.Range("H122").Value = 1872.7693
$ws.Range("I122").Value = 1805.875
$ws.Range("J122").Value = 1979.8
$ws.Range("K122").Value = 5417.625
$ws.Range("L122").Value = 5939.4
$ws.Range("M122").Value = -2967.625
$ws.Range("N122").Value = -10839.4

$ws.Range("H132").Value = 200005010
$ws.Range("I132").Value = 500005000
$ws.Range("K132").Value = 1500015000
$ws.Range("M132").Value = -1500012470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 55555920
$ws.Range("I55").Value = 71428980
$ws.Range("J55").Value = 233
$ws.Range("K55").Value = 71428980
$ws.Range("L55").Value = 233
$ws.Range("M55").Value = -71428807
$ws.Range("N55").Value = -579

$ws.Range("H122").Value = 5077.4644
$ws.Range("J122").Value = 7728.4287
$ws.Range("L122").Value = 23185.2861
$ws.Range("N122").Value = -28085.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 19232226
$ws.Range("I107").Value = 25001528
$ws.Range("J107").Value = 1214
$ws.Range("K107").Value = 75004584
$ws.Range("L107").Value = 3642
$ws.Range("M107").Value = -75002664
$ws.Range("N107").Value = -7482

$ws.Range("H132").Value = 1255
$ws.Range("I132").Value = 1265.6897
$ws.Range("J132").Value = 1100
$ws.Range("K132").Value = 3797.0691
$ws.Range("L132").Value = 1100
$ws.Range("M132").Value = -1267.0691
$ws.Range("N132").Value = -8360

$ws.Range("H136").Value = 1163.25
$ws.Range("I136").Value = 1163.25
$ws.Range("K136").Value = 3489.75
$ws.Range("M136").Value = -939.75
